$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: convert a cell that currently holds a plain number into a
# text-placeholder cell (used by this report for "no data" / "n/a" values,
# e.g. "0" or "***.*") while reusing the same cell style ("General" / s=14)
# that the rest of the placeholder cells in the row already use. We borrow
# the formatting from column A of the same row, which is always styled with
# that shared "text" style.
# ---------------------------------------------------------------------------
function Set-PlaceholderText($addr, $text, $donorRow) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range("A" + $donorRow).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# Report header: volume/issue number and the reporting week dates.
# These live inside multi-run shared strings, so we edit just the
# sub-string that changed via the Characters() API (keeps the rest of the
# rich text run untouched).
# ---------------------------------------------------------------------------

# "Volume 30   Number  17" -> "Volume 30   Number  18"
$ws.Range("A8").Characters(21, 2).Text = "18"

# "Report Covering the Week  4/24/2023  Through  4/30/2023"
#   -> "Report Covering the Week  5/1/2023  Through  5/7/2023"
$ws.Range("C9").Characters(27, 9).Text = "5/1/2023"
$ws.Range("C9").Characters(46, 9).Text = "5/7/2023"

# ---------------------------------------------------------------------------
# Weekly crime-complaint statistics table (rows 14-30)
# ---------------------------------------------------------------------------

# Row 14 - Murder: this week 2023 count now unreported ("0" placeholder)
Set-PlaceholderText "F14" "0" 14

# Row 15 - Rape
$ws.Range("N15").Value = -81.818181818181

# Row 16 - Robbery
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 3
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 57
$ws.Range("J16").Value = 59
$ws.Range("K16").Value = -3.389830508474
$ws.Range("L16").Value = 5.555555555555
$ws.Range("M16").Value = -19.718309859154
$ws.Range("N16").Value = -79.120879120879

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = -13.333333333333
$ws.Range("I17").Value = 86
$ws.Range("J17").Value = 79
$ws.Range("K17").Value = 8.860759493670
$ws.Range("L17").Value = 24.637681159420
$ws.Range("M17").Value = 36.507936507936
$ws.Range("N17").Value = -59.047619047619

# Row 18 - Burglary
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -66.666666666666
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = -57.142857142857
$ws.Range("I18").Value = 42
$ws.Range("J18").Value = 66
$ws.Range("K18").Value = -36.363636363636
$ws.Range("L18").Value = -14.285714285714
$ws.Range("M18").Value = -2.325581395348
$ws.Range("N18").Value = -86.363636363636

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 18
$ws.Range("G19").Value = 25
$ws.Range("H19").Value = -28
$ws.Range("I19").Value = 119
$ws.Range("J19").Value = 116
$ws.Range("K19").Value = 2.586206896551
$ws.Range("L19").Value = 17.821782178217
$ws.Range("M19").Value = 22.680412371134
$ws.Range("N19").Value = -16.783216783216

# Row 20 - G.L.A.
$ws.Range("C20").Value = 2
Set-PlaceholderText "D20" "0" 20
Set-PlaceholderText "E20" "***.*" 20
$ws.Range("F20").Value = 5
$ws.Range("H20").Value = 150
$ws.Range("I20").Value = 20
$ws.Range("K20").Value = 25
$ws.Range("L20").Value = 17.647058823529
$ws.Range("M20").Value = 233.333333333333
$ws.Range("N20").Value = -63.636363636363

# Row 21 - TOTAL
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = -33.333333333333
$ws.Range("F21").Value = 60
$ws.Range("G21").Value = 78
$ws.Range("H21").Value = -23.076923076923
$ws.Range("I21").Value = 329
$ws.Range("J21").Value = 339
$ws.Range("K21").Value = -2.949852507374
$ws.Range("L21").Value = 11.904761904761
$ws.Range("M21").Value = 15.438596491228
$ws.Range("N21").Value = -67.490118577075

# Row 22 - Transit
Set-PlaceholderText "D22" "0" 22
Set-PlaceholderText "E22" "***.*" 22
$ws.Range("G22").Value = 3
$ws.Range("M22").Value = -37.5

# Row 23 - Housing
$ws.Range("C23").Value = 2
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 4
$ws.Range("H23").Value = -42.857142857142
$ws.Range("I23").Value = 30
$ws.Range("J23").Value = 34
$ws.Range("K23").Value = -11.764705882352
$ws.Range("L23").Value = 7.142857142857
$ws.Range("M23").Value = 172.727272727273

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = 52.941176470588
$ws.Range("F24").Value = 103
$ws.Range("H24").Value = 33.766233766233
$ws.Range("I24").Value = 420
$ws.Range("J24").Value = 378
$ws.Range("K24").Value = 11.111111111111
$ws.Range("L24").Value = 25.373134328358
$ws.Range("M24").Value = 30.030959752322

# Row 25 - Misd. Assault
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 600
$ws.Range("F25").Value = 42
$ws.Range("G25").Value = 37
$ws.Range("H25").Value = 13.513513513513
$ws.Range("I25").Value = 147
$ws.Range("J25").Value = 157
$ws.Range("K25").Value = -6.369426751592
$ws.Range("L25").Value = -1.342281879194
$ws.Range("M25").Value = -12.5

# Row 26 - UCR Rape*
$ws.Range("L26").Value = -42.857142857142

# Row 27 - Other Sex Crimes
Set-PlaceholderText "C27" "0" 27
$ws.Range("E27").Value = -100
$ws.Range("J27").Value = 22
$ws.Range("K27").Value = -31.818181818181
$ws.Range("L27").Value = 0

# Row 28 - Shooting Vic.
Set-PlaceholderText "D28" "0" 28
Set-PlaceholderText "E28" "***.*" 28
Set-PlaceholderText "F28" "0" 28
$ws.Range("H28").Value = -100
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = -81.818181818181

# Row 29 - Shooting Inc.
Set-PlaceholderText "D29" "0" 29
Set-PlaceholderText "E29" "***.*" 29
Set-PlaceholderText "F29" "0" 29
$ws.Range("H29").Value = -100
$ws.Range("L29").Value = -16.666666666666
$ws.Range("N29").Value = -82.758620689655

# Row 30 - Hate Crimes
Set-PlaceholderText "F30" "0" 30
